# Apply "Add data for 2022-03-11" update:
#  - sheet name / label text moves from "through 03-02" to "through 03-03"
#  - March row (row 4) and Total row (row 5) get refreshed figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-03-03"

# 2. Update the March row label (shared string) to match.
$ws.Range("A4").Value = "March (through 03-03)"

# 3. Refresh the March row (row 4) counts for each year column (B:I = 2015..2022).
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 16

# 4. Refresh the Total row (row 5) to match the new March figures.
$ws.Range("B5").Value = 39
$ws.Range("C5").Value = 90
$ws.Range("D5").Value = 133
$ws.Range("E5").Value = 140
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = 145
$ws.Range("H5").Value = 351
$ws.Range("I5").Value = 316
